# Auto-generated: append new Lancers listings and refresh timestamps (2025-12-01 12:38:59 JST)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop every existing hyperlink object up front -- we re-create them below in row order
# so the rId numbering / <hyperlinks> ordering comes out exactly like a fresh export.
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2025-12-01 12:38:59"
$ws.Range("B2").Value = "【急募】生成AI×業務効率化の実装を支援するエンジニア募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5444662"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化"

# Row 3
$ws.Range("A3").Value = "2025-12-01 12:38:59"
$ws.Range("B3").Value = "製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5439165"
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = "🔥AI,Ai"

# Row 4
$ws.Range("A4").Value = "2025-12-01 12:38:59"
$ws.Range("B4").Value = "外部WEB予約サイト一元管理システム開発|長期保守パートナー募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5444378"
$ws.Range("G4").Value = 170
$ws.Range("H4").Value = "◆開発,システム開発 ◇サイト"

# Row 5
$ws.Range("A5").Value = "2025-12-01 12:38:59"
$ws.Range("B5").Value = "【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5444489"
$ws.Range("G5").Value = 155
$ws.Range("H5").Value = "◆開発,Node.js"

# Row 6
$ws.Range("A6").Value = "2025-12-01 12:38:59"
$ws.Range("B6").Value = "【Excelでのマクロ作成】リサーチツールの作成【スクレイピング】"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5445149"
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = "◆ツール,スクレイピング"

# Row 7
$ws.Range("A7").Value = "2025-12-01 12:38:59"
$ws.Range("B7").Value = "【急募】Googledriveのロール管理・共有設定の専門家募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5444395"
$ws.Range("G7").Value = 38
$ws.Range("H7").Value = "◇管理"

# Row 8
$ws.Range("A8").Value = "2025-12-01 12:38:59"
$ws.Range("B8").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "~ 5,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = "◇管理"

# Row 9
$ws.Range("A9").Value = "2025-12-01 12:38:59"
$ws.Range("B9").Value = "DAO構築。ブロックチェーンとスマートコントラクトの専門家募集"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5445105"
$ws.Range("G9").Value = 18

# Row 10
$ws.Range("A10").Value = "2025-12-01 12:38:59"
$ws.Range("B10").Value = "ホスティング業務を担当してくれる方を探しています!"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5445080"
$ws.Range("G10").Value = 18

# Row 11
$ws.Range("A11").Value = "2025-12-01 12:38:59"
$ws.Range("B11").Value = "【急募】Amazonフラットファイル(ブラウズノード検証)に詳しい方を探しています"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5444446"
$ws.Range("G11").Value = 18

# Row 12
$ws.Range("A12").Value = "2025-12-01 12:38:59"
$ws.Range("B12").Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Range("G12").Value = 13

# Row 13
$ws.Range("A13").Value = "2025-12-01 12:38:59"
$ws.Range("B13").Value = "comfyui(paperspace)でエロ動画のループ物を作成したいです。その方法を教えてください"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5444370"
$ws.Range("G13").Value = 10

# Re-create the F-column hyperlinks in row order (rId1..rId12), then normalize
# the visual style back onto the workbook's existing built-in "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5444662") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439165") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5444378") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5444489") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5445149") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5444395") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5418064") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5445105") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5445080") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5444446") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5443568") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5444370") | Out-Null

$ws.Range("F2:F13").Style = "Hyperlink"

